# Cut Sheet Express - fill in the hard coded test/demo data for the cut
# sheet grid (rows 3-4 continue the 1..N numbering started in row 2,
# row 5 adds a row of placeholder "a" values so the combined approach -
# reading real data + falling back to hard coded fill - can be compared).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nudge the "Model #:" header (C1) so it carries an (empty/general)
# alignment record, matching the rest of the header tidy-up pass.
$ws.Range("C1").HorizontalAlignment = 1

# Row 3 continues straight on from row 2 (which already holds 1-20):
# 21-40 across columns A-T.
$value = 21
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(3, $col).Value = $value
    $value++
}

# Row 4 continues the same series: 41-60 across columns A-T.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(4, $col).Value = $value
    $value++
}

# New row 5: hard coded placeholder text ("a") in every column, right
# aligned so it reads like a filled-in form field.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(5, $col).Value = "a"
}
$ws.Range("A5:T5").HorizontalAlignment = -4152

# Move the active selection past the new data, mirroring where Excel
# left the cursor after the edit.
[void]$ws.Range("A6:XFD141").Select()
